$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---------------------------
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Give the new headers the same formatting as the existing header cells
# (bold font, thin border, centered alignment) by copying the format from
# the adjacent H1 header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-assert the values (PasteSpecial only touches formatting, but make sure
# nothing was disturbed).
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows 2-33 -------------------------------------------------------
# Column I is a constant 1, column J mirrors the existing column H value.
for ($r = 2; $r -le 33; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
